$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.384.83'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.567.01'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''1.001'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = '''286.58'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '''0.3743'
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('D8').Value = '''0.3274'
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').Value = '''45.41'
$ws.Range('E9').Value = '  -5.71%  '
$ws.Range('E10').Value = '  +2.41%  '
$ws.Range('D11').Value = '''0.07417'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '''1.001'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').Value = '''20.50'
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('D14').Value = '''5.841'
$ws.Range('E14').Value = '  -2.03%  '
$ws.Range('D15').Value = '''6.837'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '1.565.27'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '''0.00001099'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '''0.06712'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '''85.87'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = '''1.000'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = '''6.349'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').Value = '''16.27'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '''11.69'
$ws.Range('E23').Value = '  -2.61%  '
$ws.Range('D24').Value = '22.390.38'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '''2.313'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').Value = '''2.553'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '''151.29'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = '''19.35'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').Value = '''4.913'
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('D30').Value = '''123.41'
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').Value = '1.743.36'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').Value = '''1.056'
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('B33').Value = 'WEMIXTOKEN'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''1.943'
$ws.Range('E33').Value = '  -2.58%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.925'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('D35').Value = '''9.631'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').Value = '''0.08262'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '''0.02387'
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '''0.06324'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D40').Value = '''0.2188'
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('D41').Value = '''5.254'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').Value = '''11.10'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('D43').Value = '''0.6101'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('D46').Value = '''3.747'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '''0.5919'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').Value = '''2.008'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').Value = '''123.81'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '''1.180'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').Value = '''0.07152'
$ws.Range('E51').Value = '  -0.78%  '

Write-Host "Applied cryptos update"